$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.04
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.04
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.03
$ws.Range("K2").Value = 950
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.26
$ws.Range("O2").Value = 1.01
$ws.Range("P2").Value = 1.26
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 1.04
$ws.Range("U2").Value = 1.04
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000
$ws.Range("G3").Value = 1.38
$ws.Range("P3").Value = 2.62
$ws.Range("R3").Value = 1.66
$ws.Range("S3").Value = 2.44
$ws.Range("U3").Value = 2.02
$ws.Range("Y3").Value = 42
$ws.Range("F4").Value = 5.5
$ws.Range("G4").Value = 5.8
$ws.Range("H4").Value = 1.69
$ws.Range("K4").Value = 4.4
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 5.1
$ws.Range("P4").Value = 2.38
$ws.Range("Q4").Value = 1.68
$ws.Range("R4").Value = 1.56
$ws.Range("S4").Value = 2.66
$ws.Range("T4").Value = 1.73
$ws.Range("U4").Value = 2.28
$ws.Range("Y4").Value = 11.5
$ws.Range("AB4").Value = 24
$ws.Range("AC4").Value = 10
$ws.Range("AE4").Value = 16.5
$ws.Range("AM4").Value = 85
$ws.Range("F5").Value = 1.53
$ws.Range("I5").Value = 2.9
$ws.Range("Q5").Value = 1.97
$ws.Range("F7").Value = 8
$ws.Range("K7").Value = 5.7
$ws.Range("R7").Value = 1.71
$ws.Range("S7").Value = 2.36
$ws.Range("AH7").Value = 22
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 3.45
$ws.Range("N8").Value = 4.7
$ws.Range("P8").Value = 2.24
$ws.Range("R8").Value = 1.48
$ws.Range("F9").Value = 1.72
$ws.Range("G9").Value = 1.74
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 5.2
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 2.58
$ws.Range("R9").Value = 1.63
$ws.Range("S9").Value = 2.48
$ws.Range("U9").Value = 2.48
$ws.Range("AB9").Value = 12.5
$ws.Range("AC9").Value = 10
$ws.Range("AG9").Value = 9.800000000000001
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 19
$ws.Range("AK9").Value = 17.5
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 7.6
$ws.Range("AO9").Value = 140
$ws.Range("G10").Value = 1.43
$ws.Range("I10").Value = 9.4
$ws.Range("J10").Value = 5.2
$ws.Range("K10").Value = 5.5
$ws.Range("N10").Value = 6
$ws.Range("F11").Value = 1.33
$ws.Range("G11").Value = 1.34
$ws.Range("H11").Value = 10
$ws.Range("I11").Value = 10.5
$ws.Range("U11").Value = 2.42
$ws.Range("AN11").Value = 3.5
$ws.Range("F12").Value = 2.4
$ws.Range("G12").Value = 2.44
$ws.Range("H12").Value = 3.2
$ws.Range("R12").Value = 1.41
$ws.Range("AC12").Value = 7.8
$ws.Range("F13").Value = 1.99
$ws.Range("G13").Value = 2.18
$ws.Range("H13").Value = 3.75
$ws.Range("I13").Value = 4.5
$ws.Range("J13").Value = 3.65
